# Weekly CompStat report refresh: new crime data collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header text: bump the report week number and the covered date range.
# ---------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/23/2024  Through  12/29/2024"

# ---------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------
$ws.Range("I14").Value = 1
$ws.Range("K14").Value = -91.666666666666
$ws.Range("L14").Value = -85.714285714285
$ws.Range("M14").Value = -92.307692307692
$ws.Range("N14").Value = -95.238095238095

# ---------------------------------------------------------------
# Row 15 - Rape (D15/E15 switch from numbers to the "0"/"***.*" placeholder text)
# ---------------------------------------------------------------
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G15").Value = 2
$ws.Range("M15").Value = -52.380952380952
$ws.Range("N15").Value = -64.285714285714

# ---------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -54.545454545454
$ws.Range("I16").Value = 125
$ws.Range("J16").Value = 163
$ws.Range("K16").Value = -23.312883435582
$ws.Range("L16").Value = 7.758620689655
$ws.Range("M16").Value = -55.830388692579
$ws.Range("N16").Value = -87.636003956478

# ---------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = 14.285714285714
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = 17.142857142857
$ws.Range("I17").Value = 477
$ws.Range("J17").Value = 469
$ws.Range("K17").Value = 1.705756929637
$ws.Range("L17").Value = 21.374045801526
$ws.Range("M17").Value = 62.244897959183
$ws.Range("N17").Value = -39.312977099236

# ---------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = -63.157894736842
$ws.Range("I18").Value = 102
$ws.Range("J18").Value = 143
$ws.Range("K18").Value = -28.671328671328
$ws.Range("L18").Value = -17.073170731707
$ws.Range("M18").Value = -67.096774193548
$ws.Range("N18").Value = -93.996468510888

# ---------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -90
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = -38.709677419354
$ws.Range("J19").Value = 415
$ws.Range("K19").Value = -12.048192771084
$ws.Range("L19").Value = 6.725146198830
$ws.Range("M19").Value = -11.622276029055
$ws.Range("N19").Value = -37.925170068027

# ---------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -58.333333333333
$ws.Range("I20").Value = 125
$ws.Range("J20").Value = 156
$ws.Range("K20").Value = -19.871794871794
$ws.Range("L20").Value = -8.088235294117
$ws.Range("M20").Value = -34.895833333333
$ws.Range("N20").Value = -90.699404761904

# ---------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = -30
$ws.Range("I21").Value = 1215
$ws.Range("J21").Value = 1378
$ws.Range("K21").Value = -11.828737300435
$ws.Range("L21").Value = 7.427055702917
$ws.Range("M21").Value = -21.460892049127
$ws.Range("N21").Value = -77.929155313351

# ---------------------------------------------------------------
# Row 23 - Transit (C23 switches from the "0" placeholder text to a number)
# ---------------------------------------------------------------
$ws.Range("D23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -37.5
$ws.Range("I23").Value = 78
$ws.Range("J23").Value = 86
$ws.Range("K23").Value = -9.302325581395
$ws.Range("L23").Value = 34.482758620689
$ws.Range("M23").Value = 56

# ---------------------------------------------------------------
# Row 24 - Housing
# ---------------------------------------------------------------
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -31.818181818181
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = -25.714285714285
$ws.Range("I24").Value = 1176
$ws.Range("J24").Value = 1304
$ws.Range("K24").Value = -9.815950920245
$ws.Range("L24").Value = -1.836393989983
$ws.Range("M24").Value = -21.442885771543

# ---------------------------------------------------------------
# Row 25 - Petit Larceny
# ---------------------------------------------------------------
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -38.461538461538
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = -11.627906976744
$ws.Range("I25").Value = 518
$ws.Range("J25").Value = 503
$ws.Range("K25").Value = 2.982107355864
$ws.Range("L25").Value = 29.5

# ---------------------------------------------------------------
# Row 26 - Retail Theft
# ---------------------------------------------------------------
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 157.142857142857
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = 7.317073170731
$ws.Range("I26").Value = 766
$ws.Range("J26").Value = 721
$ws.Range("K26").Value = 6.241331484049
$ws.Range("L26").Value = 21.011058451816
$ws.Range("M26").Value = -37.672904800650

# ---------------------------------------------------------------
# Row 27 - Misd. Assault
# ---------------------------------------------------------------
$ws.Range("D27").Value = 2
$ws.Range("G27").Value = 5
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = 16.129032258064

# ---------------------------------------------------------------
# Row 28 - UCR Rape*
# ---------------------------------------------------------------
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("I28").Value = 82
$ws.Range("J28").Value = 92
$ws.Range("K28").Value = -10.869565217391
$ws.Range("L28").Value = 6.493506493506

# ---------------------------------------------------------------
# Row 29 - Other Sex Crimes
# ---------------------------------------------------------------
$ws.Range("M29").Value = -60.714285714285
$ws.Range("N29").Value = -90.090090090090

# ---------------------------------------------------------------
# Row 30 - Shooting Vic.
# ---------------------------------------------------------------
$ws.Range("M30").Value = -56
$ws.Range("N30").Value = -88.888888888888
